# Update cryptocurrency price/volume data and reorder a few new/changed rows
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.727.01"
$ws.Range("E2").Value = "  +0.56%  "
$ws.Range("D3").Value = "1.851.14"
$ws.Range("E3").Value = "  +2.24%  "
$ws.Range("E4").Value = "  +0.31%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "227.69"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.90%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.611"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.12%  "
$ws.Range("E7").Value = "  +0.25%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "41.96"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +15.11%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.304"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +4.15%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0692"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.33%  "
$ws.Range("E11").Value = "  +3.67%  "
$ws.Range("D12").Value = "2.117.10"
$ws.Range("E12").Value = "  +2.17%  "
$ws.Range("D13").Value = "1.864.22"
$ws.Range("E13").Value = "  +2.39%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "11.37"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.71%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.71"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +6.50%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.656"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +4.17%  "
$ws.Range("D17").Value = "34.722.61"
$ws.Range("E17").Value = "  +0.64%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "68.81"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.67%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "244.47"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.73%  "
$ws.Range("D20").Value = "0.0₃0788"
$ws.Range("E20").Value = "  +1.50%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.13"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +8.04%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.79"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +16.91%  "
$ws.Range("E23").Value = "  +0.36%  "
$ws.Range("E24").Value = "  -1.60%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "172.20"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.33%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.92"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.73%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "17.89"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.60%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.123"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.16%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.31"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +6.14%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.94"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.16%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.00"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.39%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0530"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.65%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.89"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +5.02%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "91.08"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +12.32%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.667"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.00%  "
$ws.Range("D37").Value = "1.345.27"
$ws.Range("E37").Value = "  -1.29%  "
$ws.Range("B38").Value = "ARBITRUM"
$ws.Range("C38").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.02"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +8.76%  "
$ws.Range("B39").Value = "RenderToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.43"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.21%  "
$ws.Range("B40").Value = "TrustWalletToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.07"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.70%  "
$ws.Range("B41").Value = "VeChain"
$ws.Range("C41").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0193"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.85%  "
$ws.Range("B42").Value = "InjectiveProtocol"
$ws.Range("C42").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "14.89"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +10.87%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.25"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +7.54%  "
$ws.Range("B44").Value = "MXToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.84"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.18%  "
$ws.Range("B45").Value = "HuobiToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.44"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.82%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0519"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +4.17%  "
$ws.Range("D47").Value = "2.014.29"
$ws.Range("E47").Value = "  +2.14%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "6.03"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.64%  "
$ws.Range("E49").Value = "  +0.32%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "102.44"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.18%  "
$ws.Range("B51").Value = "Aptos"
$ws.Range("C51").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.29"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +5.85%  "
